# Apply "HJD and Palgrave Macmillan added" edit to the service-to-profession sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Re-point cell formatting (fonts/number-formats) so that the rows whose
#        "look" changes end up with the correct existing style.
#        Order matters: copy from rows that have not yet been overwritten.
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A8:C8").PasteSpecial(-4122) | Out-Null        # row 8 gets row 3's original look (Calibri/general)

$ws.Range("A4:C4").Copy() | Out-Null
$ws.Range("A3:C3").PasteSpecial(-4122) | Out-Null         # row 3 gets row 4's original look (Times/Arial/date)

$ws.Range("A5:C5").Copy() | Out-Null
$ws.Range("A4:C4").PasteSpecial(-4122) | Out-Null         # row 4 gets row 5's look (Times/Arial/date)

$excel.CutCopyMode = 0

# --- 2) Update the cell text/content for the rows that changed.
$ws.Range("A3").Value2 = "Board Member (Inaugural)"
$ws.Range("B3").Value2 = "2020 - Present "
$ws.Range("C3").Value2 = "The Korean Association for Public Diplomacy"

$ws.Range("A4").Value2 = "Executive Committee Member"
$ws.Range("B4").Value2 = "2020 - Present "
$ws.Range("C4").Value2 = "International Studies Association International Communication Section (ISA ICOMM)"

$ws.Range("A5").Value2 = "Section Chair"
$ws.Range("B5").Value2 = "2022 - Present"
$ws.Range("C5").Value2 = "International Studies Association International Communication Section (ISA ICOMM)"

$ws.Range("A6").Value2 = "Editor"
$ws.Range("B6").Value2 = "2022 - Present"
$ws.Range("C6").Value2 = "Korean Journal of International Studies"

$ws.Range("A7").Value2 = "Advisory Board Member"
$ws.Range("B7").Value2 = "2022 - Present"
$ws.Range("C7").Value2 = "The Hague Journal of Diplomacy"

$ws.Range("A8").Value2 = "Associate Editor"
$ws.Range("B8").Value2 = "2017 - 2022"
$ws.Range("C8").Value2 = "Journal of Contemporary Eastern Asia"

$ws.Range("A10").Value2 = "Founder and Organizer"
$ws.Range("B10").Value2 = "2019 - 2020"
$ws.Range("C10").Value2 = "Korea Public Diplomacy Colloquium (which later became the foundation of the Korean Association for Public Diplomacy)"

# --- 3) Widen the first two columns so the longer text fits.
$ws.Columns("A").ColumnWidth = 18.330729166666668
$ws.Columns("B").ColumnWidth = 15.666666666666666

# --- 4) Leave the cursor/selection on C2, like in the saved file.
$ws.Range("C2").Select() | Out-Null
